# Apply professional table formatting to the Quick Reference Card's
# "Common Scenarios" table: blue header row with white bold text,
# alternating light-gray shading on data rows, and evenly-sized columns.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Word COM RGB colors are encoded as 0xBBGGRR (blue high byte, red low
# byte) rather than the familiar 0xRRGGBB hex notation used in OOXML.
# This little helper converts a standard "RRGGBB" hex string into the
# decimal value Word expects for Shading.BackgroundPatternColor /
# Font.Color assignments.
function RGBHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$blue   = RGBHex("4472C4")
$gray   = RGBHex("E7E6E6")
$white  = RGBHex("FFFFFF")

# ---------------------------------------------------------------------
# 1) Equalize the four column widths (2030/1624/2436/1827 -> 2160 each,
#    i.e. 101.5pt/81.2pt/121.8pt/91.35pt -> 108pt) so header text no
#    longer wraps.
for ($col = 1; $col -le $t.Columns.Count; $col++) {
    $t.Columns.Item($col).Width = 108
}

# ---------------------------------------------------------------------
# 2) Header row (row 1): blue fill + white bold text on every cell.
$headerRow = 1
for ($col = 1; $col -le $t.Columns.Count; $col++) {
    $cell = $t.Cell($headerRow, $col)
    $cell.Shading.BackgroundPatternColor = $blue

    $rng = $d.Range($cell.Range.Start, $cell.Range.End - 1)
    $rng.Font.Color = $white
}

# ---------------------------------------------------------------------
# 3) Shaded data rows: rows 3 and 5 (the "Browser close" and "Fixed
#    version" rows) get light-gray shading; rows 2 and 4 stay white.
$shadedRows = @(3, 5)
foreach ($r in $shadedRows) {
    for ($col = 1; $col -le $t.Columns.Count; $col++) {
        $cell = $t.Cell($r, $col)
        $cell.Shading.BackgroundPatternColor = $gray
    }
}
